$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 59, shifting existing rows 59-88 down to 61-90
$ws.Rows.Item(59).Insert()
$ws.Rows.Item(59).Insert()

# Populate new row 59
$ws.Cells.Item(59,1).Value = 10
$ws.Cells.Item(59,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(59,3).Value = "La Araucanía"
$ws.Cells.Item(59,4).Value = 44518
$ws.Cells.Item(59,5).Value = 9
$ws.Cells.Item(59,6).Value = 100112031
$ws.Cells.Item(59,7).Value = "Poroto verde"
$ws.Cells.Item(59,8).Value = "Sin especificar"
$ws.Cells.Item(59,9).Value = "Primera"
$ws.Cells.Item(59,10).Value = 200
$ws.Cells.Item(59,11).Value = 3000
$ws.Cells.Item(59,12).Value = 3000
$ws.Cells.Item(59,13).Value = 3000
$ws.Cells.Item(59,14).Value = "`$/kilo"
$ws.Cells.Item(59,15).Value = "Región Metropolitana"
$ws.Cells.Item(59,16).Value = 3000
$ws.Cells.Item(59,17).Value = 1
$ws.Cells.Item(59,18).Value = "Hortaliza"

# Populate new row 60
$ws.Cells.Item(60,1).Value = 10
$ws.Cells.Item(60,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(60,3).Value = "La Araucanía"
$ws.Cells.Item(60,4).Value = 44518
$ws.Cells.Item(60,5).Value = 9
$ws.Cells.Item(60,6).Value = 100112031
$ws.Cells.Item(60,7).Value = "Poroto verde"
$ws.Cells.Item(60,8).Value = "Sin especificar"
$ws.Cells.Item(60,9).Value = "Primera"
$ws.Cells.Item(60,10).Value = 55
$ws.Cells.Item(60,11).Value = 70000
$ws.Cells.Item(60,12).Value = 70000
$ws.Cells.Item(60,13).Value = 70000
$ws.Cells.Item(60,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(60,15).Value = "Región Metropolitana"
$ws.Cells.Item(60,16).Value = 2800
$ws.Cells.Item(60,17).Value = 25
$ws.Cells.Item(60,18).Value = "Hortaliza"
